$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 54 values
$ws.Range("A54").Value2 = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -0.7200474048664085
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -0.2284091334091687

# Copy formatting from row 53 (A53 has the date style) to row 54 for column A
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)  # xlPasteFormats
